$d = $word.ActiveDocument

$replacements = @(
    @("17×72=", "35×25="),
    @("80×99=", "61×55="),
    @("58×98=", "35×88="),
    @("91×27=", "45×35="),
    @("40×45=", "54×35="),
    @("90×30=", "66×37="),
    @("35×73=", "58×26="),
    @("37×36=", "44×29="),
    @("82×96=", "80×93="),
    @("97×17=", "28×90="),
    @("33×36=", "76×20="),
    @("42×93=", "16×74="),
    @("69×52=", "19×43="),
    @("96×95=", "80×93="),
    @("90×32=", "30×69="),
    @("74×30=", "37×50="),
    @("17×14=", "13×27="),
    @("27×99=", "43×50="),
    @("39×92=", "36×42="),
    @("31×47=", "30×21="),
    @("90×13=", "69×84="),
    @("39×70=", "11×85="),
    @("53×28=", "81×35="),
    @("97×84=", "56×35="),
    @("47×45=", "28×33=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
